$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (bestFit-style width matching the authored column, ~15.55 chars)
$ws.Columns.Item(1).ColumnWidth = 14.6

# Update dates in B3 and B4 (43692 -> 43696, i.e. 2019-08-15 -> 2019-08-19)
$ws.Range("B3").Value = 43696
$ws.Range("B4").Value = 43696

# Update the active selection to C9
$ws.Range("C9").Select()
